# C5-PowerPoint.pptx edit
#
# 1) The table on slide 6 ("SOURCES OF FINANCE") switches its table style
#    from the deck's custom "Table_0" style to the built-in
#    {7321AC04-159A-4C1D-A5A7-CD77BA7DFC66} table style.
#
# 2) The presentation's theme colour scheme (ppt/theme/theme2.xml, the
#    theme used by the slide master / whole deck) is swapped from the
#    "Integral" palette to the default "Office Theme" palette (the
#    palette that used to live in ppt/theme/theme1.xml, which is only
#    wired up to the notes master).

$p = $ppt.ActivePresentation

# --- 1) Table style -------------------------------------------------------
# Find the (only) table in the deck -- slide 6's "SOURCES OF FINANCE" table
# -- wherever it lives, and switch it onto the built-in table style.

for ($si = 1; $si -le $p.Slides.Count; $si++) {
    $slide = $p.Slides.Item($si)
    for ($i = 1; $i -le $slide.Shapes.Count; $i++) {
        $shp = $slide.Shapes.Item($i)
        if ($shp.HasTable) {
            $shp.Table.ApplyStyle("{7321AC04-159A-4C1D-A5A7-CD77BA7DFC66}")
        }
    }
}

# --- 2) Theme colours -------------------------------------------------------
# msoThemeColorDark1=1, Light1=2, Dark2=3, Light2=4, Accent1..6=5-10,
# Hyperlink=11, FollowedHyperlink=12. RGB values use the COM BGR-packed
# "0x00BBGGRR" encoding PowerPoint's RGB property expects.

$officeThemeColors = @{
    1  = 0x000000   # Dark 1    -> 000000
    2  = 0xFFFFFF   # Light 1   -> FFFFFF
    3  = 0x6A5444   # Dark 2    -> 44546A
    4  = 0xE6E6E7   # Light 2   -> E7E6E6
    5  = 0xD59B5B   # Accent 1  -> 5B9BD5
    6  = 0x317DED   # Accent 2  -> ED7D31
    7  = 0xA5A5A5   # Accent 3  -> A5A5A5
    8  = 0x00C0FF   # Accent 4  -> FFC000
    9  = 0xC47244   # Accent 5  -> 4472C4
    10 = 0x47AD70   # Accent 6  -> 70AD47
    11 = 0xC16305   # Hyperlink -> 0563C1
    12 = 0x724F95   # FollowedHyperlink -> 954F72
}

$themeColors = $p.Slides.Item(1).ThemeColorScheme
foreach ($idx in $officeThemeColors.Keys) {
    $themeColors.Colors($idx).RGB = $officeThemeColors[$idx]
}
